# Auto-generated edit script: insert 2022-Q4 sheet + update totals sheet
$wb = $excel.ActiveWorkbook

$refSheet = $wb.Worksheets.Item(2)
$templateSheet = $wb.Worksheets.Item('2022-Q3')
$q4 = $wb.Worksheets.Add($refSheet)
$q4.Name = '2022-Q4'

# --- reuse existing header/index styling from the 2022-Q3 sheet ---
$templateSheet.Range("A2:A9").Copy()
$q4.Range("A2:A22").PasteSpecial(-4122)
$templateSheet.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

# --- populate 2022-Q4 sheet ---
$q4.Range("B2:G22").NumberFormat = "@"
$q4.Range('B1').Value = '基金代码'
$q4.Range('C1').Value = '基金名称'
$q4.Range('D1').Value = '基金规模'
$q4.Range('E1').Value = '股票总仓位'
$q4.Range('F1').Value = '仓位占比'
$q4.Range('G1').Value = '持有市值(亿元)'
$q4.Range('H1').Value = '仓位排名'
$q4.Range('A2').Value = 0
$q4.Range('B2').Value = '008903'
$q4.Range('C2').Value = '广发科技先锋混合'
$q4.Range('D2').Value = '129.95'
$q4.Range('E2').Value = '93.81'
$q4.Range('F2').Value = '4.92'
$q4.Range('G2').Value = '6.3935'
$q4.Range('H2').Value = 10
$q4.Range('A3').Value = 1
$q4.Range('B3').Value = '162703'
$q4.Range('C3').Value = '广发小盘成长混合（LOF）A'
$q4.Range('D3').Value = '96.14'
$q4.Range('E3').Value = '93.72'
$q4.Range('F3').Value = '5.95'
$q4.Range('G3').Value = '5.7203'
$q4.Range('H3').Value = 7
$q4.Range('A4').Value = 2
$q4.Range('B4').Value = '003745'
$q4.Range('C4').Value = '广发多元新兴股票'
$q4.Range('D4').Value = '34.42'
$q4.Range('E4').Value = '91.97'
$q4.Range('F4').Value = '5.28'
$q4.Range('G4').Value = '1.8174'
$q4.Range('H4').Value = 7
$q4.Range('A5').Value = 3
$q4.Range('B5').Value = '009086'
$q4.Range('C5').Value = '鹏华价值共赢两年持有期混合'
$q4.Range('D5').Value = '11.11'
$q4.Range('E5').Value = '94.37'
$q4.Range('F5').Value = '3.93'
$q4.Range('G5').Value = '0.4366'
$q4.Range('H5').Value = 9
$q4.Range('A6').Value = 4
$q4.Range('B6').Value = '009132'
$q4.Range('C6').Value = '广发小盘成长混合（LOF）C'
$q4.Range('D6').Value = '3.80'
$q4.Range('E6').Value = '93.72'
$q4.Range('F6').Value = '5.95'
$q4.Range('G6').Value = '0.2261'
$q4.Range('H6').Value = 7
$q4.Range('A7').Value = 5
$q4.Range('B7').Value = '159870'
$q4.Range('C7').Value = '鹏华中证细分化工产业主题ETF'
$q4.Range('D7').Value = '8.99'
$q4.Range('E7').Value = '98.40'
$q4.Range('F7').Value = '2.09'
$q4.Range('G7').Value = '0.1879'
$q4.Range('H7').Value = 8
$q4.Range('A8').Value = 6
$q4.Range('B8').Value = '516020'
$q4.Range('C8').Value = '华宝中证细分化工产业主题ETF'
$q4.Range('D8').Value = '3.64'
$q4.Range('E8').Value = '97.94'
$q4.Range('F8').Value = '2.08'
$q4.Range('G8').Value = '0.0757'
$q4.Range('H8').Value = 8
$q4.Range('A9').Value = 7
$q4.Range('B9').Value = '516120'
$q4.Range('C9').Value = '富国中证细分化工产业主题ETF'
$q4.Range('D9').Value = '2.94'
$q4.Range('E9').Value = '99.53'
$q4.Range('F9').Value = '2.11'
$q4.Range('G9').Value = '0.0620'
$q4.Range('H9').Value = 8
$q4.Range('A10').Value = 8
$q4.Range('B10').Value = '516220'
$q4.Range('C10').Value = '国泰中证细分化工产业主题ETF'
$q4.Range('D10').Value = '1.83'
$q4.Range('E10').Value = '98.74'
$q4.Range('F10').Value = '2.05'
$q4.Range('G10').Value = '0.0375'
$q4.Range('H10').Value = 9
$q4.Range('A11').Value = 9
$q4.Range('B11').Value = '012272'
$q4.Range('C11').Value = '渤海汇金创新价值一年持有期混合'
$q4.Range('D11').Value = '0.89'
$q4.Range('E11').Value = '85.56'
$q4.Range('F11').Value = '3.74'
$q4.Range('G11').Value = '0.0333'
$q4.Range('H11').Value = 4
$q4.Range('A12').Value = 10
$q4.Range('B12').Value = '003780'
$q4.Range('C12').Value = '鹏华兴悦定期开放灵活配置混合'
$q4.Range('D12').Value = '2.03'
$q4.Range('E12').Value = '28.03'
$q4.Range('F12').Value = '1.40'
$q4.Range('G12').Value = '0.0284'
$q4.Range('H12').Value = 4
$q4.Range('A13').Value = 11
$q4.Range('B13').Value = '003142'
$q4.Range('C13').Value = '鹏华弘达灵活配置混合A'
$q4.Range('D13').Value = '1.18'
$q4.Range('E13').Value = '39.70'
$q4.Range('F13').Value = '2.23'
$q4.Range('G13').Value = '0.0263'
$q4.Range('H13').Value = 8
$q4.Range('A14').Value = 12
$q4.Range('B14').Value = '014408'
$q4.Range('C14').Value = '创金合信兴选产业趋势一年封闭混合A'
$q4.Range('D14').Value = '1.31'
$q4.Range('E14').Value = '50.54'
$q4.Range('F14').Value = '1.40'
$q4.Range('G14').Value = '0.0183'
$q4.Range('H14').Value = 8
$q4.Range('A15').Value = 13
$q4.Range('B15').Value = '014409'
$q4.Range('C15').Value = '创金合信兴选产业趋势一年封闭混合C'
$q4.Range('D15').Value = '0.92'
$q4.Range('E15').Value = '50.54'
$q4.Range('F15').Value = '1.40'
$q4.Range('G15').Value = '0.0129'
$q4.Range('H15').Value = 8
$q4.Range('A16').Value = 14
$q4.Range('B16').Value = '516690'
$q4.Range('C16').Value = '银华中证细分化工产业主题ETF'
$q4.Range('D16').Value = '0.57'
$q4.Range('E16').Value = '97.79'
$q4.Range('F16').Value = '2.08'
$q4.Range('G16').Value = '0.0119'
$q4.Range('H16').Value = 8
$q4.Range('A17').Value = 15
$q4.Range('B17').Value = '003143'
$q4.Range('C17').Value = '鹏华弘达灵活配置混合C'
$q4.Range('D17').Value = '0.11'
$q4.Range('E17').Value = '39.70'
$q4.Range('F17').Value = '2.23'
$q4.Range('G17').Value = '0.0025'
$q4.Range('H17').Value = 8
$q4.Range('A18').Value = 16
$q4.Range('B18').Value = '013527'
$q4.Range('C18').Value = '嘉实中证细分化工产业主题指数A'
$q4.Range('D18').Value = '0.11'
$q4.Range('E18').Value = '94.81'
$q4.Range('F18').Value = '2.02'
$q4.Range('G18').Value = '0.0022'
$q4.Range('H18').Value = 8
$q4.Range('A19').Value = 17
$q4.Range('B19').Value = '015897'
$q4.Range('C19').Value = '天弘中证细分化工指数C'
$q4.Range('D19').Value = '0.10'
$q4.Range('E19').Value = '94.84'
$q4.Range('F19').Value = '2.02'
$q4.Range('G19').Value = '0.0020'
$q4.Range('H19').Value = 8
$q4.Range('A20').Value = 18
$q4.Range('B20').Value = '001474'
$q4.Range('C20').Value = '兴银丰盈灵活配置混合'
$q4.Range('D20').Value = '0.08'
$q4.Range('E20').Value = '90.33'
$q4.Range('F20').Value = '2.13'
$q4.Range('G20').Value = '0.0017'
$q4.Range('H20').Value = 10
$q4.Range('A21').Value = 19
$q4.Range('B21').Value = '015896'
$q4.Range('C21').Value = '天弘中证细分化工指数A'
$q4.Range('D21').Value = '0.05'
$q4.Range('E21').Value = '94.84'
$q4.Range('F21').Value = '2.02'
$q4.Range('G21').Value = '0.0010'
$q4.Range('H21').Value = 8
$q4.Range('A22').Value = 20
$q4.Range('B22').Value = '013528'
$q4.Range('C22').Value = '嘉实中证细分化工产业主题指数C'
$q4.Range('D22').Value = '0.02'
$q4.Range('E22').Value = '94.81'
$q4.Range('F22').Value = '2.02'
$q4.Range('G22').Value = '0.0004'
$q4.Range('H22').Value = 8

# --- rewrite 总计 (totals) sheet rows 2-7 ---
$total = $wb.Worksheets.Item(1)
$total.Range('A2').Value = 0
$total.Range('B2').Value = '2022-Q4'
$total.Range('C2').Value = 21
$total.Range('D2').Value = 15.1
$total.Range('A3').Value = 1
$total.Range('B3').Value = '2022-Q3'
$total.Range('C3').Value = 8
$total.Range('D3').Value = 12.83
$total.Range('A4').Value = 2
$total.Range('B4').Value = '2022-Q2'
$total.Range('C4').Value = 7
$total.Range('D4').Value = 27.51
$total.Range('A5').Value = 3
$total.Range('B5').Value = '2022-Q1'
$total.Range('C5').Value = 12
$total.Range('D5').Value = 39.6
$total.Range('A6').Value = 4
$total.Range('B6').Value = '2021-Q4'
$total.Range('C6').Value = 14
$total.Range('D6').Value = 50.42
$total.Range('A7').Value = 5
$total.Range('B7').Value = '2021-Q3'
$total.Range('C7').Value = 14
$total.Range('D7').Value = 54.73
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)
